$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 6 (Site ID 5) with new "Karasjok Camping" site data
$ws.Range("B6").Value = "Karasjok Camping"
$ws.Range("C6").Value = 25.487036705017001
$ws.Range("D6").Value = 69.467986997867399
$ws.Range("E6").Value = 440783.85183459503
$ws.Range("F6").Value = 7707278.4545117402

# Update the active selection on the sheet
$ws.Range("B18").Select()
